$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.451.71'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').Value = '1.791.95'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '226.53'
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = '0.558'
$ws.Range('E6').Value = '  +1.88%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '32.77'
$ws.Range('E8').Value = '  +3.31%  '
$ws.Range('D9').Value = '0.297'
$ws.Range('E9').Value = '  +1.89%  '
$ws.Range('D10').Value = '0.0693'
$ws.Range('E10').Value = '  +0.61%  '
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').Value = '2.049.49'
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('D13').Value = '11.09'
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').Value = '1.792.15'
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').Value = '0.636'
$ws.Range('E15').Value = '  +1.97%  '
$ws.Range('D16').Value = '34.416.42'
$ws.Range('E16').Value = '  +0.84%  '
$ws.Range('D17').Value = '4.27'
$ws.Range('E17').Value = '  +2.14%  '
$ws.Range('D18').Value = '68.77'
$ws.Range('E18').Value = '  +0.77%  '
$ws.Range('D19').Value = '247.33'
$ws.Range('E19').Value = '  -0.10%  '
$ws.Range('D20').Value = '0.0₃0799'
$ws.Range('E20').Value = '  +2.87%  '
$ws.Range('D21').Value = '11.30'
$ws.Range('E21').Value = '  +3.67%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('E23').Value = '  +1.15%  '
$ws.Range('D24').Value = '2.06'
$ws.Range('E24').Value = '  +0.87%  '
$ws.Range('D25').Value = '164.77'
$ws.Range('E25').Value = '  +2.33%  '
$ws.Range('D26').Value = '7.25'
$ws.Range('E26').Value = '  +1.02%  '
$ws.Range('D27').Value = '16.52'
$ws.Range('E27').Value = '  +1.20%  '
$ws.Range('E28').Value = '  +2.39%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '3.92'
$ws.Range('E30').Value = '  +7.80%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '3.80'
$ws.Range('E31').Value = '  +3.23%  '
$ws.Range('E32').Value = '  +0.44%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.23'
$ws.Range('E33').Value = '  -0.28%  '
$ws.Range('E34').Value = '  +1.35%  '
$ws.Range('D35').Value = '1.417.70'
$ws.Range('E35').Value = '  -2.05%  '
$ws.Range('D36').Value = '2.58'
$ws.Range('E36').Value = '  +5.46%  '
$ws.Range('D37').Value = '0.672'
$ws.Range('E37').Value = '  +2.74%  '
$ws.Range('E38').Value = '  +0.42%  '
$ws.Range('E39').Value = '  +1.63%  '
$ws.Range('D40').Value = '84.99'
$ws.Range('E40').Value = '  +5.59%  '
$ws.Range('E41').Value = '  +0.86%  '
$ws.Range('D42').Value = '0.935'
$ws.Range('E42').Value = '  +1.40%  '
$ws.Range('E43').Value = '  +2.29%  '
$ws.Range('D44').Value = '13.53'
$ws.Range('E44').Value = '  +0.67%  '
$ws.Range('D45').Value = '0.0523'
$ws.Range('E45').Value = '  +2.77%  '
$ws.Range('D46').Value = '6.05'
$ws.Range('E46').Value = '  -0.23%  '
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').Value = '1.949.18'
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('D49').Value = '105.47'
$ws.Range('E49').Value = '  -0.17%  '
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').Value = '0.0₆0128'
$ws.Range('E51').Value = '  -5.58%  '
